# Omaha_Cal_Info_GS05MOAS-GL003_00001.xlsx
#
# Commit: "Updated Global Glider Cal and Ingest sheets"
#   - Changed Cal scattering angle to 140
#   - Changed angular resolution to 1.13
#   - (data-source "recovered" -> "recovered_host" change belongs to a
#     different workbook in the same commit; not present in this file)
#
# The Asset_Cal_Info sheet holds the calibration-coefficient rows:
#   row 2 -> CC_scattering_angle      (column F)
#   row 4 -> CC_angular_resolution    (column F)

$wb = $excel.ActiveWorkbook

$calSheet = $wb.Worksheets.Item("Asset_Cal_Info")

# Update the two calibration coefficient values.
$calSheet.Range("F2").Value = 140
$calSheet.Range("F4").Value = 1.13

# The author ended the session on the Asset_Cal_Info tab with F20
# selected (moving focus away from the Moorings tab).
$calSheet.Activate() | Out-Null
$calSheet.Range("F20").Select() | Out-Null
